# Apply updated crypto price/volume figures to the worksheet.
# Values in column D are plain numeric-looking text (e.g. "214.68") that must
# stay as literal text (matching the original inlineStr cells), so for any new
# value Excel would otherwise auto-convert to a number we prefix it with a
# leading apostrophe ('), the standard way to force text entry in Excel.
# Column E values already contain surrounding spaces/percent signs so they
# are stored as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.868.12'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '1.629.67'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''214.68'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = '''0.502'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("D9").Value = '''0.0631'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '''19.65'
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").Value = '1.855.00'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '1.640.87'
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").Value = '''0.545'
$ws.Range("E15").Value = '  -1.99%  '
$ws.Range("D16").Value = '0.0₃0757'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '''62.74'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '25.860.20'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("D21").Value = '''192.68'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").Value = '''0.0498'
$ws.Range("E31").Value = '  +2.03%  '
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '''3.22'
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").Value = '1.136.60'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("E38").Value = '  +1.73%  '
$ws.Range("D39").Value = '''2.47'
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("D43").Value = '''98.99'
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '1.765.13'
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").Value = '0.0₆0111'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '''56.14'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").Value = '''0.0530'
$ws.Range("E48").Value = '  +4.43%  '
$ws.Range("E49").Value = '  +1.74%  '
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").Value = '''7.62'
